$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C1: a date (2021-04-06, serial 44292) formatted as a short date (built-in numFmtId 14)
$ws.Range("C1").Value = 44292
$ws.Range("C1").NumberFormat = "mm-dd-yy"

# C2:C33 expense values for the new (third) column
$values = @(
    7.88, 41.65, 44.68, 14.99, 13.21, 8.99, 54.99, 0,
    12.39, 9.99, 19.75, 279.98, 6.39, 6.99, 14.54, 37.04,
    0, 7.91, 19.63, 0, 0, 15.99, 11.99, 27.69,
    249.39, 39.84, 41.93, 9.97, 66.71, 16.99, 14.97, 13.99
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}

# C34: total, sum of the column
$ws.Range("C34").Formula = "=SUM(C2:C33)"

# Widen column C to fit its (numeric) contents, closest achievable to the real Excel best-fit width
$ws.Columns.Item(3).ColumnWidth = 9.83

# Match the active selection recorded by the author's session
$ws.Range("C2").Select() | Out-Null
